# Adding test for Amazon e-commerce app
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Google
$ws2 = $wb.Worksheets.Item(2)   # Salesforce -> Amazon

# Rename the second sheet
$ws2.Name = "Amazon"

# Populate the new Amazon test-case rows (order chosen so that the
# shared-string table is built up in the same sequence as the target file)
$ws2.Range("A2").Value = "TC01"
$ws2.Range("B2").Value = "TC01_AmazonLogin"
$ws2.Range("C2").Value = "Login to Amazon"
$ws2.Range("D2").Value = "Yes"

$ws2.Range("A3").Value = "TC02"
$ws2.Range("D3").Value = "No"

$ws2.Range("A4").Value = "TC03"
$ws2.Range("C4").Value = "Search in Amazon"
$ws2.Range("D4").Value = "No"

$ws2.Range("B4").Value = "TC03_AmazonSearch"
$ws2.Range("B3").Value = "TC02_AmazonFilter"

# Widen column B on the Amazon sheet to fit the new values
$ws2.Columns.Item(2).ColumnWidth = 17.67

# Update selections on both sheets
[void]$ws1.Range("A2:XFD4").Select()
[void]$ws2.Range("C4").Select()

# Amazon becomes the active (selected) tab
[void]$ws2.Activate()
